$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (CVR, Year, TCV, Løsning, Opsagt dato, Quarter) added from
# the newest export (7th October 2024).
#
# Column A (CVR) holds numeric-looking values but must be stored as TEXT
# (shared string), exactly like the existing rows 2-14. Briefly applying a
# text number format forces Excel to keep the typed value as a string, then
# clearing the format returns the cell to the default (unstyled) look that
# the rest of column A already uses.
#
# Values are written column-by-column (both new rows at a time) so that new
# shared-string entries land in the same grouping order as the source data.

# --- Column A (CVR) --------------------------------------------------------
$ws.Range("A15:A16").NumberFormat = "@"
$ws.Range("A15").Value = "33638760"
$ws.Range("A16").Value = "38538071"
$ws.Range("A15:A16").ClearFormats()

# --- Column B (Year) --------------------------------------------------------
$ws.Range("B15").Value = 2024
$ws.Range("B16").Value = 2024

# --- Column C (Beløb 12 mdr. (TCV)) -----------------------------------------
$ws.Range("C15").Value = 127477
$ws.Range("C16").Value = 129216

# --- Column D (Løsning) ------------------------------------------------------
$ws.Range("D15").Value = "Løn/HR og Time"
$ws.Range("D16").Value = "EasyCruit"

# --- Column E (Opsagt dato) --------------------------------------------------
# Reuse the workbook's existing custom date format so it shares the same
# style index as the other date cells instead of creating a brand-new style.
$ws.Range("E15:E16").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E15").Value = "2024-01-04"
$ws.Range("E16").Value = "2024-06-07"

# --- Column H (Quarter) ------------------------------------------------------
$ws.Range("H15").Value = "2024Q1"
$ws.Range("H16").Value = "2024Q2"

# --- Column I (TCV_range) ----------------------------------------------------
$ws.Range("I15").Value = "120000-140000"
$ws.Range("I16").Value = "120000-140000"
